$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the numeric value in A20 with the text "gym"
$ws.Range("A20").Value = "gym"

# Add a new row 21 with "bring groceries" (same text as used in row 8 / A8)
$ws.Range("A21").Value = "bring groceries"

# Update the view state: scroll so row 4 is the top-left visible row,
# and select cell B21
$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Range("B21").Select()
